# Complete project upload to GitHub
# - Insert a new "vvvv" user as row 2 (pushing the existing "ss" and "vasu"
#   rows down by one), and append a new "aa" user as the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 ("ss"/21), shifting the existing
# data rows (21/ss -> row 3, 5/vasu -> row 4) down by one.
$ws.Rows("2:2").Insert()

# Populate the newly inserted row 2 with the "vvvv" signup.
$ws.Range("A2").Value = 30
$ws.Range("B2").Value = "vvvv"
$ws.Range("C2").Value = "vvvv@gmail.com"
$ws.Range("D2").Value = "student"
$ws.Range("E2").Value = "20/1/2026"
$ws.Range("F2").Value = "6:28:59 pm"

# Append the new "aa" signup as row 5.
$ws.Range("A5").Value = 33
$ws.Range("B5").Value = "aa"
$ws.Range("C5").Value = "aa@gmail.com"
$ws.Range("D5").Value = "student"
$ws.Range("E5").Value = "27/1/2026"
$ws.Range("F5").Value = "11:31:50 am"
